$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so values like "27.455.40" or "22.10"
# are not auto-coerced into numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.455.40'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.570.10'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '207.38'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '22.10'
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '1.795.84'
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('D13').Value = '1.588.11'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('E15').Value = '  -3.08%  '
$ws.Range('D16').Value = '63.17'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '27.450.76'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '213.62'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '4.12'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').Value = '9.68'
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('D25').Value = '152.95'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('D26').Value = '6.86'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').Value = '3.20'
$ws.Range('E32').Value = '  -1.89%  '
$ws.Range('D33').Value = '1.359.12'
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('D39').Value = '0.531'
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('E40').Value = '  +1.31%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '0.973'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = '64.13'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('D47').Value = '1.707.15'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('D48').Value = '85.34'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '0.0₇0996'
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').Value = '0.0956'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').Value = '0.0494'
$ws.Range('E51').Value = '  -0.70%  '

# Restore the default (Normal) style on column D so the text-format
# override above does not leave a stray style index on the cells.
$ws.Range("D2:D51").Style = "Normal"
